$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.82"
$ws.Range("E2").Value = "'0.69%"
$ws.Range("D3").Value = "'41.01"
$ws.Range("E3").Value = "'2.11%"
$ws.Range("D4").Value = "'5.131"
$ws.Range("E4").Value = "'1.70%"
$ws.Range("D5").Value = "'0.07629"
$ws.Range("E5").Value = "'0.44%"
$ws.Range("D6").Value = "'4.277"
$ws.Range("E6").Value = "'0.40%"
$ws.Range("D7").Value = "'1.624"
$ws.Range("E7").Value = "'1.90%"
$ws.Range("D8").Value = "'2.459"
$ws.Range("E8").Value = "'0.00%"
$ws.Range("D9").Value = "'0.9094"
$ws.Range("E9").Value = "'0.18%"
$ws.Range("D10").Value = "'0.1170"
$ws.Range("E10").Value = "'17.04%"
$ws.Range("D11").Value = "'0.1798"
$ws.Range("E11").Value = "'2.31%"
$ws.Range("D12").Value = "'0.09147"
$ws.Range("E12").Value = "'1.48%"
$ws.Range("D13").Value = "'0.04253"
$ws.Range("E13").Value = "'-2.77%"
$ws.Range("D14").Value = "'0.1044"
$ws.Range("E14").Value = "'-0.72%"
$ws.Range("D15").Value = "'0.001259"
$ws.Range("E15").Value = "'2.02%"
$ws.Range("D16").Value = "'0.005856"
$ws.Range("E16").Value = "'0.61%"
$ws.Range("D17").Value = "'3.354"
$ws.Range("E17").Value = "'-0.32%"
$ws.Range("E18").Value = "'-0.71%"
$ws.Range("D19").Value = "'6.892"
$ws.Range("E19").Value = "'0.16%"
$ws.Range("D20").Value = "'0.1394"
$ws.Range("E20").Value = "'2.45%"
$ws.Range("D21").Value = "'0.2706"
$ws.Range("D22").Value = "'0.04031"
$ws.Range("E22").Value = "'-3.21%"
$ws.Range("E23").Value = "'4.59%"
$ws.Range("D24").Value = "'0.004080"
$ws.Range("E24").Value = "'0.53%"
$ws.Range("E25").Value = "'-2.35%"
$ws.Range("D26").Value = "'0.0003751"
$ws.Range("E38").Value = "'0.40%"
$ws.Range("D39").Value = "'0.05230"
$ws.Range("E39").Value = "'1.78%"
$ws.Range("D40").Value = "'0.007798"
$ws.Range("E40").Value = "'-0.66%"
$ws.Range("D41").Value = "'0.1302"
$ws.Range("E41").Value = "'0.02%"
$ws.Range("D42").Value = "'0.006802"
$ws.Range("E42").Value = "'-4.05%"
$ws.Range("D43").Value = "'0.001934"
$ws.Range("E43").Value = "'-1.95%"
$ws.Range("D44").Value = "'0.008069"
$ws.Range("E44").Value = "'-3.47%"
$ws.Range("E45").Value = "'-7.75%"
$ws.Range("D46").Value = "'0.00006898"
$ws.Range("E46").Value = "'6.89%"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("D48").Value = "'0.09456"
$ws.Range("E48").Value = "'1,858.46%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'-0.05%"
